# Insert a new data row for "Vega Modelo de Temuco - Repollo" right after
# the header/first-data rows, at row 799, pushing the existing rows
# 799-854 down to 800-855 (dimension grows from A1:R854 to A1:R855).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(799).Insert()

$ws.Range("A799").Value = 10
$ws.Range("B799").Value = "Vega Modelo de Temuco"
$ws.Range("C799").Value = "La Araucanía"
$ws.Range("D799").Value = "2023-01-05"
$ws.Range("E799").Value = 9
$ws.Range("F799").Value = 100112006
$ws.Range("G799").Value = "Repollo"
$ws.Range("H799").Value = "Crespo record"
$ws.Range("I799").Value = "Primera"
$ws.Range("J799").Value = 300
$ws.Range("K799").Value = 1200
$ws.Range("L799").Value = 1200
$ws.Range("M799").Value = 1200
$ws.Range("N799").Value = "$/unidad"
$ws.Range("O799").Value = "Provincia de Cautín"
$ws.Range("P799").Value = 1200
$ws.Range("Q799").Value = 1
$ws.Range("R799").Value = "Hortaliza"
